$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the existing "99 - Other Identifier" row (old
# row 10). That row then slides down to row 11 untouched (same values,
# same style), and the freshly inserted row 10 is filled in with the new
# "09 - EDUMIS" entry.
$ws.Rows("10").Insert()

$numericLookingRange = $ws.Range("A10:B10")
$numericLookingRange.NumberFormat = "@"

$ws.Range("A10").Value = "1.0"
$ws.Range("B10").Value = "09"
$ws.Range("C10").Value = "EDUMIS"
$ws.Range("D10").Value = "Identifier for educational providers registered with the Ministry of Education's Education Management Information System (EDUMIS)"
